# ============================================================
# "study in office 221010" - 2022-10-10 study session
#  1) Refresh the "time" (col E) timestamp on cards reviewed today
#  2) Append newly studied cards (rows 122-149, part = 10월_Nick_Drama Topic 11)
# ============================================================

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Update review timestamps on existing rows ---
$ws.Cells.Item(3, 5).Value = 44844.76965959491
$ws.Cells.Item(4, 5).Value = 44844.83561008102
$ws.Cells.Item(5, 5).Value = 44844.83495538194
$ws.Cells.Item(10, 5).Value = 44844.77333452546
$ws.Cells.Item(11, 5).Value = 44844.77070377315
$ws.Cells.Item(13, 5).Value = 44844.7803833912
$ws.Cells.Item(17, 5).Value = 44844.77354201389
$ws.Cells.Item(18, 5).Value = 44844.77345179398
$ws.Cells.Item(23, 5).Value = 44844.82119972222
$ws.Cells.Item(25, 5).Value = 44844.78163292824
$ws.Cells.Item(27, 5).Value = 44844.80394146991
$ws.Cells.Item(31, 5).Value = 44844.82685549768
$ws.Cells.Item(34, 5).Value = 44830.77858189815
$ws.Cells.Item(38, 5).Value = 44844.82576273148
$ws.Cells.Item(43, 5).Value = 44844.82852371528
$ws.Cells.Item(44, 5).Value = 44844.82908050926
$ws.Cells.Item(45, 5).Value = 44844.83249679398
$ws.Cells.Item(49, 5).Value = 44844.77856605324
$ws.Cells.Item(50, 5).Value = 44844.77778479167
$ws.Cells.Item(51, 5).Value = 44844.82916106482
$ws.Cells.Item(54, 5).Value = 44844.83259344907
$ws.Cells.Item(55, 5).Value = 44844.76793635417
$ws.Cells.Item(57, 5).Value = 44844.82869263889
$ws.Cells.Item(58, 5).Value = 44844.79723451389
$ws.Cells.Item(60, 5).Value = 44844.81070938658
$ws.Cells.Item(61, 5).Value = 44844.79402392361
$ws.Cells.Item(62, 5).Value = 44844.79416642361
$ws.Cells.Item(64, 5).Value = 44844.77262538194
$ws.Cells.Item(66, 5).Value = 44844.80398449074
$ws.Cells.Item(67, 5).Value = 44844.77116351852
$ws.Cells.Item(69, 5).Value = 44844.83264717593
$ws.Cells.Item(71, 5).Value = 44844.80362780092
$ws.Cells.Item(72, 5).Value = 44844.83548678241
$ws.Cells.Item(73, 5).Value = 44844.82676074074
$ws.Cells.Item(83, 5).Value = 44844.80407037037
$ws.Cells.Item(84, 5).Value = 44844.77745317129
$ws.Cells.Item(85, 5).Value = 44844.79814195602
$ws.Cells.Item(86, 5).Value = 44844.79907678241
$ws.Cells.Item(87, 5).Value = 44844.83490446759
$ws.Cells.Item(88, 5).Value = 44844.77707355324
$ws.Cells.Item(89, 5).Value = 44844.77985939815
$ws.Cells.Item(92, 5).Value = 44844.77016797454
$ws.Cells.Item(94, 5).Value = 44844.80621981482
$ws.Cells.Item(96, 5).Value = 44844.77455603009
$ws.Cells.Item(97, 5).Value = 44844.80108430555
$ws.Cells.Item(98, 5).Value = 44844.83616904904
$ws.Cells.Item(99, 5).Value = 44844.79886037037
$ws.Cells.Item(100, 5).Value = 44844.77286246528
$ws.Cells.Item(101, 5).Value = 44844.82894560185
$ws.Cells.Item(102, 5).Value = 44844.83463002315
$ws.Cells.Item(106, 5).Value = 44844.83054753472
$ws.Cells.Item(108, 5).Value = 44844.83522643519
$ws.Cells.Item(109, 5).Value = 44844.83475425926
$ws.Cells.Item(110, 5).Value = 44844.83307723379
$ws.Cells.Item(111, 5).Value = 44844.83503821759

# --- 2) Append new rows 122-149 ---

# Copy column-A number style (bold, bordered, centered) down onto the new rows
$ws.Range("A121").Copy()
$ws.Range("A122:A149").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A122").Value = 120
$ws.Range("B122").Value = "그렇게 해서 여자애들이 결국 죽는다`n[그게 여자애들이 결국 죽는 방법이야]"
$ws.Range("C122").Value = "That's how girls end up dead"
$ws.Range("D122").Value = 0
$ws.Range("F122").Value = "10월_Nick_Drama Topic 11"

$ws.Range("A123").Value = 121
$ws.Range("B123").Value = "외국인이랑 우연히 말을 했는데 말을 못 했거든. 그때 배우기로 결심했지`n[요전 날 우연히 외국인과 대화를 나누었는데, 아무 말도 할 수가 없었어요. 그때 나는 영어를 배우기로 결심했다.]"
$ws.Range("C123").Value = "The other day, I happened to have a conversation with a foreginer, but I wasmn't able say anything at all. That's when I decided to learn English."
$ws.Range("D123").Value = 0
$ws.Range("F123").Value = "10월_Nick_Drama Topic 11"

$ws.Range("A124").Value = 122
$ws.Range("B124").Value = "학원에 오면 공부를 하지 말고 말을 더 해 독해 말고 그래야 회화가 늘지`n[그냥 읽는 게 아니라 학원이 있는 곳에서 영어를 더 많이 해야 할 것 같아. 그래야 영어 실력이 향상될 거야.]"
$ws.Range("C124").Value = "I think you should try to speak english more where you're at the academy, not just read. That's how your english will imporve."
$ws.Range("D124").Value = 0
$ws.Range("F124").Value = "10월_Nick_Drama Topic 11"

$ws.Range("A125").Value = 123
$ws.Range("B125").Value = "이번달에는 저축 좀 하나 했는데 뭐 또 다 써버렸지`n[이번 달에 돈을 좀 모을 수 있을 줄 알았는데, 결국 돈을 다 써버렸어.]"
$ws.Range("C125").Value = "I thought I could save some money this month, but I ended up spending all my money"
$ws.Range("D125").Value = 0
$ws.Range("F125").Value = "10월_Nick_Drama Topic 11"

$ws.Range("A126").Value = 124
$ws.Range("B126").Value = "널 좀 봐, 너 완전 홀딱 젖었어. 그것은 정말 형편없다, 너가 그런사람이다 (누구) 이것을 하기를 원했다`n[널 좀 봐. 흠뻑 젖었군요. 그건 너무 시시해! 네가 하고 싶었던 거야]"
$ws.Range("C126").Value = "look at you. You're soaked. That is so lame! You're the on who wanted to do it"
$ws.Range("D126").Value = 0
$ws.Range("F126").Value = "10월_Nick_Drama Topic 11"

$ws.Range("A127").Value = 125
$ws.Range("B127").Value = "비 쫄딱 맞았네, 다 젖엇다니까. 뉴스 봤으면 우산 챙겼을 텐데`n[비를 맞았어요. 흠뻑 젖었어요. 내가 그 뉴스를 봤더라면. 난 우산을 가져갔을 거야]"
$ws.Range("C127").Value = "I got rained on. I got totally soaked. If I had seen the news. I would've taken an umbrella"
$ws.Range("D127").Value = 0
$ws.Range("F127").Value = "10월_Nick_Drama Topic 11"

$ws.Range("A128").Value = 126
$ws.Range("B128").Value = "매일 퇴근하고 운동 갈 거라고 한 사람은 너였던거 같은데`n[퇴근 후에 운동하러 갈 사람은 너였어.]"
$ws.Range("C128").Value = "You were the one who was going to go to work out after work"
$ws.Range("D128").Value = 0
$ws.Range("F128").Value = "10월_Nick_Drama Topic 11"

$ws.Range("A129").Value = 127
$ws.Range("B129").Value = "동생이 게임 중독이 됐는데, 혼낼 수 없어. 내가 게임기를 사줬거든`n[내 남동생은 게임을 좋아했지만, 나는 그것에 대해 그에게 소리지를 수 없다. 그를 위해 그것을 산 사람은 나였다.]"
$ws.Range("C129").Value = "My brother has gotten into playing games, but I can't yell at him about it. I was the one who bought it for him"
$ws.Range("D129").Value = 0
$ws.Range("F129").Value = "10월_Nick_Drama Topic 11"

$ws.Range("A130").Value = 128
$ws.Range("B130").Value = "우리 동네에 있는 소방관들은 아주 섹시한 것으로 소문났다. 내가 그것에 분노하냐고?`n[우리 마을의 소방관들은 덥기로 유명하다. 내가 그걸 원망해?]"
$ws.Range("C130").Value = "the firemen in our town have a reputation for being hot. Do I resent that?"
$ws.Range("D130").Value = 0
$ws.Range("F130").Value = "10월_Nick_Drama Topic 11"

$ws.Range("A131").Value = 129
$ws.Range("B131").Value = "그 사람 평판이 그렇게 좋지는 않더라고 . 너무 믿지는 마 그사람이 하는 말`n[그 사람 평판이 생각보다 좋지 않아요. 그를 너무 믿지 마세요.]"
$ws.Range("C131").Value = "You know, his reputation isn't as good as I thought. Don't belive him too much."
$ws.Range("D131").Value = 0
$ws.Range("F131").Value = "10월_Nick_Drama Topic 11"

$ws.Range("A132").Value = 130
$ws.Range("B132").Value = "너 오늘따라 왜 친절 한 거야? 너 나한테 잘못한거 있는 거 맞구나`n[왜 나한테 그렇게 잘해 주는 거야? 야, 너 나한테 잘못한 게 틀림없어!]"
$ws.Range("C132").Value = "Why are you being so nice to me? Hey, you must've done something wrong to me!"
$ws.Range("D132").Value = 0
$ws.Range("F132").Value = "10월_Nick_Drama Topic 11"

$ws.Range("A133").Value = 131
$ws.Range("B133").Value = "왜 나만 매일 야근을 해야 하는지 화가 난다니까?`n[매일 밤 늦게까지 일한다는 사실이 원망스러울 수밖에 없다.]"
$ws.Range("C133").Value = "I can't help but resent the fact that I've been working late every night."
$ws.Range("D133").Value = 0
$ws.Range("F133").Value = "10월_Nick_Drama Topic 11"

$ws.Range("A134").Value = 132
$ws.Range("B134").Value = "그렇게 해야 ~하다`n[그렇게]"
$ws.Range("C134").Value = "That's how"
$ws.Range("D134").Value = 0
$ws.Range("F134").Value = "10월_Nick_Drama Topic 11"

$ws.Range("A135").Value = 133
$ws.Range("B135").Value = "결국 ~이 되다`n[결국 ~하게 되다]"
$ws.Range("C135").Value = "End up"
$ws.Range("D135").Value = 0
$ws.Range("F135").Value = "10월_Nick_Drama Topic 11"

$ws.Range("A136").Value = 134
$ws.Range("B136").Value = "죽은, 완전, 진짜로`n[죽었어]"
$ws.Range("C136").Value = "Dead"
$ws.Range("D136").Value = 0
$ws.Range("F136").Value = "10월_Nick_Drama Topic 11"

$ws.Range("A137").Value = 135
$ws.Range("B137").Value = "보다 Vs 찾다`n[보기 vs 찾기]"
$ws.Range("C137").Value = "Look at Vs Look for"
$ws.Range("D137").Value = 0
$ws.Range("F137").Value = "10월_Nick_Drama Topic 11"

$ws.Range("A138").Value = 136
$ws.Range("B138").Value = "습하다, 촉촉하다 젖었다, 홀딱젖엇다`n[습기/습기/습기/습기/습기]"
$ws.Range("C138").Value = "humid / Moisture / Wet / Soak"
$ws.Range("D138").Value = 0
$ws.Range("F138").Value = "10월_Nick_Drama Topic 11"

$ws.Range("A139").Value = 137
$ws.Range("B139").Value = "형편없는, 구리다`n[살집이 좋은]"
$ws.Range("C139").Value = "leam"
$ws.Range("D139").Value = 0
$ws.Range("F139").Value = "10월_Nick_Drama Topic 11"

$ws.Range("A140").Value = 138
$ws.Range("B140").Value = "사람을 강조하는것 너`n[바로 너야]"
$ws.Range("C140").Value = "You're the one"
$ws.Range("D140").Value = 0
$ws.Range("F140").Value = "10월_Nick_Drama Topic 11"

$ws.Range("A141").Value = 139
$ws.Range("B141").Value = "평판, 이미지`n[명성.]"
$ws.Range("C141").Value = "Reputation"
$ws.Range("D141").Value = 0
$ws.Range("F141").Value = "10월_Nick_Drama Topic 11"

$ws.Range("A142").Value = 140
$ws.Range("B142").Value = "분노를 하다, 부들부들`n[원망]"
$ws.Range("C142").Value = "Resent"
$ws.Range("D142").Value = 0
$ws.Range("F142").Value = "10월_Nick_Drama Topic 11"

$ws.Range("A143").Value = 141
$ws.Range("B143").Value = "저번에`n[요전날]"
$ws.Range("C143").Value = "The other day"
$ws.Range("D143").Value = 0
$ws.Range("F143").Value = "10월_Nick_Drama Topic 11"

$ws.Range("A144").Value = 142
$ws.Range("B144").Value = "비맞았다`n[비를 맞았다.]"
$ws.Range("C144").Value = "I got rained on"
$ws.Range("D144").Value = 0
$ws.Range("F144").Value = "10월_Nick_Drama Topic 11"

$ws.Range("A145").Value = 143
$ws.Range("B145").Value = "화내다`n[야단법석을 떨다]"
$ws.Range("C145").Value = "yell it"
$ws.Range("D145").Value = 0
$ws.Range("F145").Value = "10월_Nick_Drama Topic 11"

$ws.Range("A146").Value = 144
$ws.Range("B146").Value = "나 혼났다 선생님께`n[선생님한테 혼났어요.]"
$ws.Range("C146").Value = "I got yelled at by my teacher"
$ws.Range("D146").Value = 0
$ws.Range("F146").Value = "10월_Nick_Drama Topic 11"

$ws.Range("A147").Value = 145
$ws.Range("B147").Value = "그 상태인것 VS 변한것`n[맞으면 맞으면 맞바꾸다]"
$ws.Range("C147").Value = "be Vs get"
$ws.Range("D147").Value = 0
$ws.Range("F147").Value = "10월_Nick_Drama Topic 11"

$ws.Range("A148").Value = 146
$ws.Range("B148").Value = "존재한다고 믿는다`n[을 믿다]"
$ws.Range("C148").Value = "believe in"
$ws.Range("D148").Value = 0
$ws.Range("F148").Value = "10월_Nick_Drama Topic 11"

$ws.Range("A149").Value = 147
$ws.Range("B149").Value = "잘될거라고 믿는다`n[난 널 믿어]"
$ws.Range("C149").Value = "I believe in you"
$ws.Range("D149").Value = 0
$ws.Range("F149").Value = "10월_Nick_Drama Topic 11"

Write-Output "Applied 221010 study-session edits: updated 55 timestamps, appended 28 new rows (A122:F149)."
